$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "16+36="
$t.Cell(1,2).Range.Text = "63-15="
$t.Cell(1,3).Range.Text = "59+19="
$t.Cell(1,4).Range.Text = "55+9="
$t.Cell(1,5).Range.Text = "15+48="
$t.Cell(2,1).Range.Text = "30-11="
$t.Cell(2,2).Range.Text = "12-8="
$t.Cell(2,3).Range.Text = "27+15="
$t.Cell(2,4).Range.Text = "7+5="
$t.Cell(2,5).Range.Text = "53-36="
$t.Cell(3,1).Range.Text = "85-68="
$t.Cell(3,2).Range.Text = "71-16="
$t.Cell(3,3).Range.Text = "62-8="
$t.Cell(3,4).Range.Text = "48+27="
$t.Cell(3,5).Range.Text = "64+27="
$t.Cell(4,1).Range.Text = "93-45="
$t.Cell(4,2).Range.Text = "25-19="
$t.Cell(4,3).Range.Text = "27-9="
$t.Cell(4,4).Range.Text = "83-45="
$t.Cell(4,5).Range.Text = "70-39="
$t.Cell(5,1).Range.Text = "8+67="
$t.Cell(5,2).Range.Text = "39+8="
$t.Cell(5,3).Range.Text = "29+13="
$t.Cell(5,4).Range.Text = "21-13="
$t.Cell(5,5).Range.Text = "56+29="
$t.Cell(6,1).Range.Text = "50-32="
$t.Cell(6,2).Range.Text = "75+18="
$t.Cell(6,3).Range.Text = "47+15="
$t.Cell(6,4).Range.Text = "17+19="
$t.Cell(6,5).Range.Text = "25+18="
$t.Cell(7,1).Range.Text = "56-29="
$t.Cell(7,2).Range.Text = "92-37="
$t.Cell(7,3).Range.Text = "56-9="
$t.Cell(7,4).Range.Text = "84-35="
$t.Cell(7,5).Range.Text = "20-14="
$t.Cell(8,1).Range.Text = "9+12="
$t.Cell(8,2).Range.Text = "41-38="
$t.Cell(8,3).Range.Text = "83-37="
$t.Cell(8,4).Range.Text = "58+17="
$t.Cell(8,5).Range.Text = "26+6="
$t.Cell(9,1).Range.Text = "30-23="
$t.Cell(9,2).Range.Text = "86+8="
$t.Cell(9,3).Range.Text = "64-28="
$t.Cell(9,4).Range.Text = "49+48="
$t.Cell(9,5).Range.Text = "64+29="
$t.Cell(10,1).Range.Text = "57+39="
$t.Cell(10,2).Range.Text = "53+19="
$t.Cell(10,3).Range.Text = "62-29="
$t.Cell(10,4).Range.Text = "9+17="
$t.Cell(10,5).Range.Text = "13+39="
$t.Cell(11,1).Range.Text = "17+79="
$t.Cell(11,2).Range.Text = "72-38="
$t.Cell(11,3).Range.Text = "49+46="
$t.Cell(11,4).Range.Text = "97-59="
$t.Cell(11,5).Range.Text = "36+39="
$t.Cell(12,1).Range.Text = "72-37="
$t.Cell(12,2).Range.Text = "46+37="
$t.Cell(12,3).Range.Text = "36+37="
$t.Cell(12,4).Range.Text = "92-73="
$t.Cell(12,5).Range.Text = "71-17="
$t.Cell(13,1).Range.Text = "77+16="
$t.Cell(13,2).Range.Text = "53+38="
$t.Cell(13,3).Range.Text = "94-86="
$t.Cell(13,4).Range.Text = "66+26="
$t.Cell(13,5).Range.Text = "51-42="
$t.Cell(14,1).Range.Text = "61-35="
$t.Cell(14,2).Range.Text = "87+7="
$t.Cell(14,3).Range.Text = "65-8="
$t.Cell(14,4).Range.Text = "17+15="
$t.Cell(14,5).Range.Text = "73-28="
$t.Cell(15,1).Range.Text = "73-48="
$t.Cell(15,2).Range.Text = "72-4="
$t.Cell(15,3).Range.Text = "39+39="
$t.Cell(15,4).Range.Text = "8+69="
$t.Cell(15,5).Range.Text = "4+49="
$t.Cell(16,1).Range.Text = "81-69="
$t.Cell(16,2).Range.Text = "84-17="
$t.Cell(16,3).Range.Text = "2+69="
$t.Cell(16,4).Range.Text = "18+74="
$t.Cell(16,5).Range.Text = "62-37="
$t.Cell(17,1).Range.Text = "54+28="
$t.Cell(17,2).Range.Text = "76-67="
$t.Cell(17,3).Range.Text = "46+17="
$t.Cell(17,4).Range.Text = "28+24="
$t.Cell(17,5).Range.Text = "39+52="
$t.Cell(18,1).Range.Text = "91-86="
$t.Cell(18,2).Range.Text = "82-74="
$t.Cell(18,3).Range.Text = "41-4="
$t.Cell(18,4).Range.Text = "67+15="
$t.Cell(18,5).Range.Text = "82-64="
$t.Cell(19,1).Range.Text = "73-66="
$t.Cell(19,2).Range.Text = "72-16="
$t.Cell(19,3).Range.Text = "82-45="
$t.Cell(19,4).Range.Text = "66-39="
$t.Cell(19,5).Range.Text = "56+8="
$t.Cell(20,1).Range.Text = "47+34="
$t.Cell(20,2).Range.Text = "66-18="
$t.Cell(20,3).Range.Text = "86-7="
$t.Cell(20,4).Range.Text = "90-84="
$t.Cell(20,5).Range.Text = "60-15="